# feat: add 2022-Q1 data
#
# 1) Duplicate the "2021-Q4" sheet (same layout/format) to create "2022-Q1",
#    positioned right before "总计", then patch in the 2022-Q1 figures.
# 2) Insert a new summary row into "总计" for 2022-Q1, pushing the
#    existing rows down and renumbering their index column.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$src.Copy($total)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

function Set-TextValue($range, $text) {
    $range.Formula = "'" + $text
    $range.Style = "Normal"
}

# Row 2: 001481 / 华宝油气(QDII)美元
Set-TextValue $q1.Range("D2") "39.80"
Set-TextValue $q1.Range("E2") "94.60"
Set-TextValue $q1.Range("F2") "2.57"
Set-TextValue $q1.Range("G2") "1.0229"
$q1.Range("H2").Value = 2

# Row 3: 162411 / 华宝油气(QDII)人民币A
Set-TextValue $q1.Range("D3") "39.80"
Set-TextValue $q1.Range("E3") "94.60"
Set-TextValue $q1.Range("F3") "2.57"
Set-TextValue $q1.Range("G3") "1.0229"
$q1.Range("H3").Value = 2

# Row 4: 007844 / 华宝油气(QDII)人民币C
Set-TextValue $q1.Range("D4") "12.98"
Set-TextValue $q1.Range("E4") "94.60"
Set-TextValue $q1.Range("F4") "2.57"
Set-TextValue $q1.Range("G4") "0.3336"
$q1.Range("H4").Value = 2

# --- Update "总计" with a new leading row for 2022-Q1 ---
$ws = $wb.Worksheets.Item("总计")

$ws.Range("A2:D2").EntireRow.Insert()
$ws.Range("A2:D2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 2.38

$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2

# Restore the originally-active tab/selection (sheet edits above shift focus)
$wb.Worksheets.Item("2021-Q1").Activate()
